# Weekly Fruta/Hortaliza update: a new record for the latest week is
# prepended to the data table, pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 22 (the first data row), shifting every
# existing record (rows 22..132) down by one (to rows 23..133).
$ws.Rows("22:22").Insert()

# Populate the newly inserted row with the new weekly record. The
# non-numeric / descriptive columns mirror the record that used to sit at
# row 22 (now shifted to row 23): market id, market name, region,
# category, variety, quality, unit of sale, origin, and classification
# stay the same; only the date and the volume/price figures change.
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value = 44561
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 100112031
$ws.Range("G22").Value = "Poroto verde"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 40000
$ws.Range("L22").Value = 40000
$ws.Range("M22").Value = 40000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 1600
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
